$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0153508771929825
$ws.Range("C2").Value = 0.959429824561403
$ws.Range("D2").Value = 0.0109649122807018
$ws.Range("E2").Value = 0.766447368421053
$ws.Range("F2").Value = 0.025219298245614
$ws.Range("G2").Value = 0.00657894736842105
$ws.Range("H2").Value = 0.0219298245614035
$ws.Range("I2").Value = 0.860745614035088
$ws.Range("J2").Value = 0.00767543859649123
$ws.Range("K2").Value = 0.00328947368421053
$ws.Range("L2").Value = 0.0328947368421053
$ws.Range("M2").Value = 0.893640350877193
$ws.Range("N2").Value = 0.00109649122807018
$ws.Range("O2").Value = 0.00219298245614035
$ws.Range("P2").Value = 0.00328947368421053
$ws.Range("Q2").Value = 0.00548245614035088
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0.00219298245614035
$ws.Range("T2").Value = 0.00219298245614035
$ws.Range("U2").Value = 0.859649122807018
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0.00657894736842105
$ws.Range("X2").Value = 0.00548245614035088
$ws.Range("B3").Value = 0.00109649122807018
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.0274122807017544
$ws.Range("E3").Value = 0.0208333333333333
$ws.Range("F3").Value = 0.155701754385965
$ws.Range("G3").Value = 0.967105263157895
$ws.Range("H3").Value = 0.978070175438597
$ws.Range("I3").Value = 0.106359649122807
$ws.Range("J3").Value = 0.923245614035088
$ws.Range("K3").Value = 0.00986842105263158
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0.00109649122807018
$ws.Range("N3").Value = 0.866228070175439
$ws.Range("O3").Value = 0.966008771929825
$ws.Range("P3").Value = 0.0317982456140351
$ws.Range("Q3").Value = 0.00109649122807018
$ws.Range("R3").Value = 0.018640350877193
$ws.Range("S3").Value = 0.0043859649122807
$ws.Range("T3").Value = 0.00109649122807018
$ws.Range("U3").Value = 0.0208333333333333
$ws.Range("V3").Value = 0.00657894736842105
$ws.Range("W3").Value = 0.025219298245614
$ws.Range("X3").Value = 0.00219298245614035
$ws.Range("B4").Value = 0.975877192982456
$ws.Range("C4").Value = 0.0350877192982456
$ws.Range("D4").Value = 0.00109649122807018
$ws.Range("E4").Value = 0.0241228070175439
$ws.Range("F4").Value = 0.0307017543859649
$ws.Range("G4").Value = 0.0208333333333333
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.0043859649122807
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.00219298245614035
$ws.Range("L4").Value = 0.957236842105263
$ws.Range("M4").Value = 0.00328947368421053
$ws.Range("N4").Value = 0.121710526315789
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0.00548245614035088
$ws.Range("R4").Value = 0.0230263157894737
$ws.Range("S4").Value = 0.990131578947368
$ws.Range("T4").Value = 0.989035087719298
$ws.Range("U4").Value = 0.110745614035088
$ws.Range("V4").Value = 0.0241228070175439
$ws.Range("W4").Value = 0.964912280701754
$ws.Range("X4").Value = 0.964912280701754
$ws.Range("B5").Value = 0.00767543859649123
$ws.Range("C5").Value = 0.00548245614035088
$ws.Range("D5").Value = 0.960526315789474
$ws.Range("E5").Value = 0.18859649122807
$ws.Range("F5").Value = 0.788377192982456
$ws.Range("G5").Value = 0.0043859649122807
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.0285087719298246
$ws.Range("J5").Value = 0.0690789473684211
$ws.Range("K5").Value = 0.984649122807018
$ws.Range("L5").Value = 0.00986842105263158
$ws.Range("M5").Value = 0.101973684210526
$ws.Range("N5").Value = 0.0109649122807018
$ws.Range("O5").Value = 0.0317982456140351
$ws.Range("P5").Value = 0.964912280701754
$ws.Range("Q5").Value = 0.987938596491228
$ws.Range("R5").Value = 0.958333333333333
$ws.Range("S5").Value = 0.00328947368421053
$ws.Range("T5").Value = 0.00767543859649123
$ws.Range("U5").Value = 0.0087719298245614
$ws.Range("V5").Value = 0.969298245614035
$ws.Range("W5").Value = 0.00219298245614035
$ws.Range("X5").Value = 0.025219298245614
